$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.145.02'
$ws.Range("E2").Value = '  -2.17%  '
$ws.Range("D3").Value = '1.577.36'
$ws.Range("E3").Value = '  -1.45%  '
$ws.Range("E4").Value = '  -0.33%  '
$ws.Range("D5").Value = "'" + '209.34'
$ws.Range("E5").Value = '  -1.06%  '
$ws.Range("E6").Value = '  -3.22%  '
$ws.Range("E7").Value = '  -0.34%  '
$ws.Range("D8").Value = "'" + '0.245'
$ws.Range("E8").Value = '  -0.63%  '
$ws.Range("E9").Value = '  -1.57%  '
$ws.Range("E10").Value = '  -0.66%  '
$ws.Range("D11").Value = "'" + '0.0843'
$ws.Range("E11").Value = '  -0.46%  '
$ws.Range("D12").Value = '1.799.53'
$ws.Range("E12").Value = '  -1.43%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = "'" + '4.04'
$ws.Range("E13").Value = '  -0.12%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.573.34'
$ws.Range("E14").Value = '  -1.66%  '
$ws.Range("D15").Value = "'" + '0.513'
$ws.Range("E15").Value = '  -1.95%  '
$ws.Range("D16").Value = "'" + '64.37'
$ws.Range("E16").Value = '  -1.08%  '
$ws.Range("D17").Value = '26.149.37'
$ws.Range("E17").Value = '  -2.03%  '
$ws.Range("D18").Value = '0.0₃0728'
$ws.Range("E18").Value = '  -1.41%  '
$ws.Range("E19").Value = '  +1.59%  '
$ws.Range("D21").Value = "'" + '207.40'
$ws.Range("E21").Value = '  -1.24%  '
$ws.Range("D22").Value = "'" + '4.25'
$ws.Range("E22").Value = '  -0.93%  '
$ws.Range("E23").Value = '  -1.94%  '
$ws.Range("D24").Value = "'" + '8.88'
$ws.Range("E24").Value = '  -1.22%  '
$ws.Range("D25").Value = "'" + '143.96'
$ws.Range("E25").Value = '  +0.32%  '
$ws.Range("E26").Value = '  -0.20%  '
$ws.Range("E27").Value = '  -1.36%  '
$ws.Range("E28").Value = '  -1.93%  '
$ws.Range("E29").Value = '  -0.60%  '
$ws.Range("D30").Value = "'" + '0.0507'
$ws.Range("E30").Value = '  -0.49%  '
$ws.Range("E31").Value = '  -1.08%  '
$ws.Range("E32").Value = '  -1.85%  '
$ws.Range("E33").Value = '  +0.46%  '
$ws.Range("D34").Value = '1.280.04'
$ws.Range("E34").Value = '  -0.65%  '
$ws.Range("E35").Value = '  -0.73%  '
$ws.Range("D36").Value = "'" + '0.608'
$ws.Range("E36").Value = '  +0.98%  '
$ws.Range("E37").Value = '  -0.87%  '
$ws.Range("E38").Value = '  -5.36%  '
$ws.Range("E39").Value = '  -2.87%  '
$ws.Range("D40").Value = "'" + '0.816'
$ws.Range("E40").Value = '  -1.91%  '
$ws.Range("D41").Value = "'" + '5.56'
$ws.Range("E41").Value = '  +2.82%  '
$ws.Range("E42").Value = '  -2.29%  '
$ws.Range("D43").Value = "'" + '62.47'
$ws.Range("E43").Value = '  -0.57%  '
$ws.Range("D44").Value = "'" + '0.763'
$ws.Range("E44").Value = '  -2.45%  '
$ws.Range("D45").Value = '1.712.38'
$ws.Range("E45").Value = '  -1.40%  '
$ws.Range("D46").Value = "'" + '89.05'
$ws.Range("E46").Value = '  -1.63%  '
$ws.Range("E47").Value = '  -0.36%  '
$ws.Range("D48").Value = '0.0₆0105'
$ws.Range("E48").Value = '  -0.93%  '
$ws.Range("E49").Value = '  -1.20%  '
$ws.Range("E50").Value = '  -1.92%  '
$ws.Range("E51").Value = '  -0.20%  '
